$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; F=1; G=24.871494;          H=74.61448200000001;  I=0.5523409677800469; J=0.552340967780047;  K=3; L=1; M=2.172908;          N=6.518724;           O=0.6836640069168457; P=0.6836640069168457; Q=54.04346828455201; R=486.391214560968;   S=0.3776156392168353;  T=0.3776156392168353 }
    3  = @{ E=3; F=1; G=24.871494;          H=74.61448200000001;  I=0.5523409677800469; J=0.552340967780047;  K=3; L=1; M=0.6730553333333334; N=2.019166;           O=0.2117640075251322; P=0.2117640075251322; Q=16.73989168466801; R=150.659025162012;   S=0.1169659368574127;  T=0.1169659368574127 }
    4  = @{ E=3; F=1; G=24.871494;          H=74.61448200000001;  I=0.5523409677800469; J=0.552340967780047;  K=3; L=1; M=0.332364;            N=0.9970920000000001; O=0.104571985558022;  P=0.104571985558022;  Q=8.266389231816001; R=74.39750308634402;  S=0.05775939170579898; T=0.05775939170579899 }
    5  = @{ E=3; F=1; G=6.727847;           H=20.183541;           I=0.1494106274056591; J=0.1494106274056591; K=3; L=1; M=2.172908;          N=6.518724;           O=0.6836640069168457; P=0.6836640069168457; Q=14.618992569076;    R=131.570933121684;   S=0.1021466682081128;  T=0.1021466682081128 }
    6  = @{ E=3; F=1; G=6.727847;           H=20.183541;           I=0.1494106274056591; J=0.1494106274056591; K=3; L=1; M=0.6730553333333334; N=2.019166;           O=0.2117640075251322; P=0.2117640075251322; Q=4.528213305200667; R=40.753919746806;    S=0.03163979322626671; T=0.03163979322626671 }
    7  = @{ E=3; F=1; G=6.727847;           H=20.183541;           I=0.1494106274056591; J=0.1494106274056591; K=3; L=1; M=0.332364;            N=0.9970920000000001; O=0.104571985558022;  P=0.104571985558022;  Q=2.236094140308;    R=20.124847262772;    S=0.01562416597127959; T=0.01562416597127959 }
    8  = @{ E=3; F=1; G=13.429899;          H=40.289697;            I=0.298248404814294;  J=0.298248404814294;  K=3; L=1; M=2.172908;          N=6.518724;           O=0.6836640069168457; P=0.6836640069168457; Q=29.181934976292;   R=262.637414786628;   S=0.2039016994918977;  T=0.2039016994918977 }
    9  = @{ E=3; F=1; G=13.429899;          H=40.289697;            I=0.298248404814294;  J=0.298248404814294;  K=3; L=1; M=0.6730553333333334; N=2.019166;           O=0.2117640075251322; P=0.2117640075251322; Q=9.039065148078002; R=81.35158633270201;  S=0.06315827744145282; T=0.06315827744145283 }
    10 = @{ E=3; F=1; G=13.429899;          H=40.289697;            I=0.298248404814294;  J=0.298248404814294;  K=3; L=1; M=0.332364;            N=0.9970920000000001; O=0.104571985558022;  P=0.104571985558022;  Q=4.463614951236001; R=40.17253456112401;  S=0.03118842788094346; T=0.03118842788094346 }
}

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}
